# Scheduled data-refresh update: overwrite market-price-derived cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with refreshed values.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 111593.78
$ws.Range("I15").Value = 111593.78
$ws.Range("K15").Value = 334781.34
$ws.Range("M15").Value = -334612.34
$ws.Range("H106").Value = 58289.445
$ws.Range("I106").Value = 65013.125
$ws.Range("K106").Value = 65013.125
$ws.Range("M106").Value = -64382.125
$ws.Range("H129").Value = 1878.644
$ws.Range("J129").Value = 2143
$ws.Range("L129").Value = 6429
$ws.Range("N129").Value = -16429
$ws.Range("H132").Value = 3292296.8
$ws.Range("I132").Value = 2817.508
$ws.Range("J132").Value = 19233620
$ws.Range("K132").Value = 8452.523999999999
$ws.Range("L132").Value = 57700860
$ws.Range("M132").Value = -5922.523999999999
$ws.Range("N132").Value = -57705920
$ws.Range("H134").Value = 58034.4
$ws.Range("J134").Value = 58034.4
$ws.Range("L134").Value = 58034.4
$ws.Range("N134").Value = -68174.39999999999
$ws.Range("H138").Value = 9526244
$ws.Range("I138").Value = 18520812
$ws.Range("J138").Value = 2585.2942
$ws.Range("K138").Value = 55562436
$ws.Range("L138").Value = 7755.882599999999
$ws.Range("M138").Value = -55557296
$ws.Range("N138").Value = -18035.8826
$ws.Range("H140").Value = 86250
$ws.Range("J140").Value = 86250
$ws.Range("L140").Value = 86250
$ws.Range("N140").Value = -96610

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7673.35
$ws.Range("I32").Value = 7610.7974
$ws.Range("K32").Value = 7610.7974
$ws.Range("M32").Value = -7323.7974
$ws.Range("H61").Value = 5556530.5
$ws.Range("I61").Value = 6494474.5
$ws.Range("J61").Value = 1016.53845
$ws.Range("K61").Value = 6494474.5
$ws.Range("L61").Value = 1016.53845
$ws.Range("M61").Value = -6494262.5
$ws.Range("N61").Value = -1440.53845
$ws.Range("H132").Value = 3522536.8
$ws.Range("I132").Value = 5683304
$ws.Range("J132").Value = 1286.2963
$ws.Range("K132").Value = 17049912
$ws.Range("L132").Value = 3858.8889
$ws.Range("M132").Value = -17047382
$ws.Range("N132").Value = -8918.8889
$ws.Range("H136").Value = 5556530.5
$ws.Range("I136").Value = 6494474.5
$ws.Range("J136").Value = 1016.53845
$ws.Range("K136").Value = 19483423.5
$ws.Range("L136").Value = 3049.61535
$ws.Range("M136").Value = -19480873.5
$ws.Range("N136").Value = -8149.61535

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2162.0793
$ws.Range("I134").Value = 1193.8605
$ws.Range("J134").Value = 4243.75
$ws.Range("K134").Value = 3581.5815
$ws.Range("L134").Value = 12731.25
$ws.Range("M134").Value = -1046.5815
$ws.Range("N134").Value = -17801.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1183.3
$ws.Range("I58").Value = 508.73077
$ws.Range("K58").Value = 508.73077
$ws.Range("M58").Value = -305.73077
$ws.Range("H132").Value = 6411412
$ws.Range("I132").Value = 7043268.5
$ws.Range("J132").Value = 2582.5715
$ws.Range("K132").Value = 21129805.5
$ws.Range("L132").Value = 7747.7145
$ws.Range("M132").Value = -21127275.5
$ws.Range("N132").Value = -12807.7145
$ws.Range("H134").Value = 259717.05
$ws.Range("I134").Value = 924.4868
$ws.Range("J134").Value = 1488981.8
$ws.Range("K134").Value = 2773.4604
$ws.Range("L134").Value = 4466945.4
$ws.Range("M134").Value = -238.4603999999999
$ws.Range("N134").Value = -4472015.4
$ws.Range("H136").Value = 1183.3
$ws.Range("I136").Value = 508.73077
$ws.Range("K136").Value = 1526.19231
$ws.Range("M136").Value = 1023.80769
$ws.Range("H140").Value = 39459.8
$ws.Range("J140").Value = 39459.8
$ws.Range("L140").Value = 39459.8
$ws.Range("N140").Value = -49819.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3967.0303
$ws.Range("I129").Value = 5316.364
$ws.Range("J129").Value = 3292.3635
$ws.Range("K129").Value = 15949.092
$ws.Range("L129").Value = 9877.0905
$ws.Range("M129").Value = -10949.092
$ws.Range("N129").Value = -19877.0905
$ws.Range("H131").Value = 1693.9246
$ws.Range("I131").Value = 2338.5
$ws.Range("J131").Value = 1303.2727
$ws.Range("K131").Value = 7015.5
$ws.Range("L131").Value = 3909.8181
$ws.Range("M131").Value = -1975.5
$ws.Range("N131").Value = -13989.8181

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 59650
$ws.Range("J138").Value = 59650
$ws.Range("L138").Value = 59650
$ws.Range("N138").Value = -69930

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1610.5555
$ws.Range("I61").Value = 1459
$ws.Range("J61").Value = 1800
$ws.Range("K61").Value = 1459
$ws.Range("L61").Value = 1800
$ws.Range("M61").Value = -1257
$ws.Range("N61").Value = -2204
$ws.Range("H93").Value = 887.17645
$ws.Range("I93").Value = 846.5454999999999
$ws.Range("J93").Value = 961.6667
$ws.Range("K93").Value = 846.5454999999999
$ws.Range("L93").Value = 961.6667
$ws.Range("M93").Value = 401.4545000000001
$ws.Range("N93").Value = -3457.6667
$ws.Range("H113").Value = 1610.5555
$ws.Range("I113").Value = 1459
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1459
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 711
$ws.Range("N113").Value = -6140
$ws.Range("H132").Value = 7358841.5
$ws.Range("I132").Value = 3289.9592
$ws.Range("J132").Value = 26328422
$ws.Range("K132").Value = 9869.8776
$ws.Range("L132").Value = 78985266
$ws.Range("M132").Value = -7339.8776
$ws.Range("N132").Value = -78990326
$ws.Range("H139").Value = 57960
$ws.Range("J139").Value = 57960
$ws.Range("L139").Value = 57960
$ws.Range("N139").Value = -68240

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("H54").Value = 16416.666
$ws.Range("H106").Value = 37700
$ws.Range("J106").Value = 37700
$ws.Range("L106").Value = 37700
$ws.Range("N106").Value = -40224
$ws.Range("H113").Value = 1903.8125
$ws.Range("I113").Value = 1483.8889
$ws.Range("J113").Value = 2443.7144
$ws.Range("K113").Value = 4451.6667
$ws.Range("L113").Value = 7331.1432
$ws.Range("M113").Value = -2281.6667
$ws.Range("N113").Value = -11671.1432
$ws.Range("H122").Value = 2089.7307
$ws.Range("I122").Value = 2312
$ws.Range("J122").Value = 1486.4286
$ws.Range("K122").Value = 6936
$ws.Range("L122").Value = 4459.2858
$ws.Range("M122").Value = -4486
$ws.Range("N122").Value = -9359.2858
$ws.Range("H136").Value = 966.7954999999999
$ws.Range("I136").Value = 891.4194
$ws.Range("J136").Value = 1146.5385
$ws.Range("K136").Value = 2674.2582
$ws.Range("L136").Value = 3439.6155
$ws.Range("M136").Value = -124.2582000000002
$ws.Range("N136").Value = -8539.6155
$ws.Range("M2").ClearContents()

